$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: call_id 54 -> 74, refresh timestamps on notes/tasks
$ws.Range("A2").Value = 74
$ws.Range("M2").Value = "`n[2025-09-25 18:56:26] No Notes available."
$ws.Range("N2").Value = "`n[2025-09-25 18:56:26] No tasks found for this call."

# Row 3: call_id 55 -> 75, update notes, clear tasks
$ws.Range("A3").Value = 75
$ws.Range("M3").Value = "`n[2025-09-25 18:58:24] The user is dealing with dust allergies and wants to schedule a meeting, but the preferred time provided is outside business hours."
$ws.Range("N3").ClearContents()
